$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 424.25
$ws.Range("I31").Value = 424.25
$ws.Range("K31").Value = 1272.75
$ws.Range("M31").Value = -1042.75

$ws.Range("H64").Value = 48472.637
$ws.Range("I64").Value = 145428.58
$ws.Range("J64").Value = 3226.5334
$ws.Range("K64").Value = 145428.58
$ws.Range("L64").Value = 3226.5334
$ws.Range("M64").Value = -145180.58
$ws.Range("N64").Value = -3722.5334

$ws.Range("H67").Value = 48472.637
$ws.Range("I67").Value = 145428.58
$ws.Range("J67").Value = 3226.5334
$ws.Range("K67").Value = 145428.58
$ws.Range("L67").Value = 3226.5334
$ws.Range("M67").Value = -144570.58
$ws.Range("N67").Value = -4942.5334

$ws.Range("H129").Value = 259993.56
$ws.Range("I129").Value = 14807.429
$ws.Range("J129").Value = 296510.62
$ws.Range("K129").Value = 44422.287
$ws.Range("L129").Value = 889531.86
$ws.Range("M129").Value = -39422.287
$ws.Range("N129").Value = -899531.86

$ws.Range("H137").Value = 1473.2727
$ws.Range("I137").Value = 1200.6207
$ws.Range("K137").Value = 3601.8621
$ws.Range("M137").Value = -1051.8621

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

$ws.Range("H14").Value = 2133.3333
$ws.Range("J14").Value = 2133.3333
$ws.Range("L14").Value = 2133.3333
$ws.Range("N14").Value = -2483.3333

$ws.Range("H32").Value = 21671.412
$ws.Range("I32").Value = 3964.757
$ws.Range("K32").Value = 3964.757
$ws.Range("M32").Value = -3677.757

$ws.Range("H41").Value = 7383.7144
$ws.Range("I41").Value = 2937.3333
$ws.Range("K41").Value = 2937.3333
$ws.Range("M41").Value = -2523.3333

$ws.Range("H45").Value = 48974.81
$ws.Range("I45").Value = 100980
$ws.Range("J45").Value = 1697.3636
$ws.Range("K45").Value = 100980
$ws.Range("L45").Value = 1697.3636
$ws.Range("M45").Value = -100603
$ws.Range("N45").Value = -2451.3636

$ws.Range("H61").Value = 1961.8572
$ws.Range("I61").Value = 1359.92
$ws.Range("J61").Value = 2847.0588
$ws.Range("K61").Value = 1359.92
$ws.Range("L61").Value = 2847.0588
$ws.Range("M61").Value = -1147.92
$ws.Range("N61").Value = -3271.0588

$ws.Range("H102").Value = 65242.312
$ws.Range("I102").Value = 168819.83
$ws.Range("K102").Value = 168819.83
$ws.Range("M102").Value = -167197.83

$ws.Range("H110").Value = 50101000
$ws.Range("I110").Value = 77077780
$ws.Range("J110").Value = 1271.5714
$ws.Range("K110").Value = 77077780
$ws.Range("L110").Value = 1271.5714
$ws.Range("M110").Value = -77075735
$ws.Range("N110").Value = -5361.5714

$ws.Range("H136").Value = 1961.8572
$ws.Range("I136").Value = 1359.92
$ws.Range("J136").Value = 2847.0588
$ws.Range("K136").Value = 4079.76
$ws.Range("L136").Value = 8541.1764
$ws.Range("M136").Value = -1529.76
$ws.Range("N136").Value = -13641.1764

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 113790.11
$ws.Range("I105").Value = 86014.164
$ws.Range("J105").Value = 169342
$ws.Range("K105").Value = 86014.164
$ws.Range("L105").Value = 169342
$ws.Range("M105").Value = -84267.164
$ws.Range("N105").Value = -172836

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28111.732
$ws.Range("I31").Value = 1342.1428
$ws.Range("J31").Value = 42526.13
$ws.Range("K31").Value = 1342.1428
$ws.Range("L31").Value = 42526.13
$ws.Range("M31").Value = -1047.1428
$ws.Range("N31").Value = -43116.13

$ws.Range("H32").Value = 19936.666
$ws.Range("J32").Value = 19900
$ws.Range("L32").Value = 19900
$ws.Range("N32").Value = -20532

$ws.Range("H33").Value = 6647.4287
$ws.Range("I33").Value = 5406.4
$ws.Range("J33").Value = 9750
$ws.Range("K33").Value = 5406.4
$ws.Range("L33").Value = 9750
$ws.Range("M33").Value = -5027.4
$ws.Range("N33").Value = -10508

$ws.Range("H34").Value = 28111.732
$ws.Range("I34").Value = 1342.1428
$ws.Range("J34").Value = 42526.13
$ws.Range("K34").Value = 1342.1428
$ws.Range("L34").Value = 42526.13
$ws.Range("M34").Value = -1140.1428
$ws.Range("N34").Value = -42930.13

$ws.Range("H47").Value = 27450
$ws.Range("J47").Value = 29933.334
$ws.Range("L47").Value = 29933.334
$ws.Range("N47").Value = -31065.334

$ws.Range("H48").Value = 18000
$ws.Range("J48").Value = 18000
$ws.Range("L48").Value = 18000
$ws.Range("N48").Value = -18952

$ws.Range("H55").Value = 9871.625
$ws.Range("J55").Value = 15875
$ws.Range("L55").Value = 15875
$ws.Range("N55").Value = -16505

$ws.Range("H82").Value = 35000
$ws.Range("J82").Value = 35000
$ws.Range("L82").Value = 35000
$ws.Range("N82").Value = -35722

$ws.Range("H85").Value = 35000
$ws.Range("J85").Value = 35000
$ws.Range("L85").Value = 35000
$ws.Range("N85").Value = -37496

$ws.Range("H94").Value = 867.3333
$ws.Range("I94").Value = 753
$ws.Range("J94").Value = 894.2353000000001
$ws.Range("K94").Value = 753
$ws.Range("L94").Value = 894.2353000000001
$ws.Range("M94").Value = -302
$ws.Range("N94").Value = -1796.2353

$ws.Range("H105").Value = 2608.5715
$ws.Range("J105").Value = 2084.4
$ws.Range("L105").Value = 2084.4
$ws.Range("N105").Value = -5578.4

$ws.Range("H133").Value = 62666.332
$ws.Range("J133").Value = 62666.332
$ws.Range("L133").Value = 62666.332
$ws.Range("N133").Value = -67726.33199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 111111420
$ws.Range("I4").Value = 121.42857
$ws.Range("J4").Value = 500001000
$ws.Range("K4").Value = 364.28571
$ws.Range("L4").Value = 1500003000
$ws.Range("M4").Value = -252.28571
$ws.Range("N4").Value = -1500003224

$ws.Range("H121").Value = 1473193.1
$ws.Range("J121").Value = 2678169.2
$ws.Range("L121").Value = 8034507.600000001
$ws.Range("N121").Value = -8037127.600000001

$ws.Range("H131").Value = 814.55206
$ws.Range("I131").Value = 417
$ws.Range("J131").Value = 860.77905
$ws.Range("K131").Value = 1251
$ws.Range("L131").Value = 2582.33715
$ws.Range("M131").Value = 3789
$ws.Range("N131").Value = -12662.33715

$ws.Range("H133").Value = 6480
$ws.Range("I133").Value = 700
$ws.Range("J133").Value = 10333.333
$ws.Range("K133").Value = 2100
$ws.Range("L133").Value = 30999.999
$ws.Range("M133").Value = 2960
$ws.Range("N133").Value = -41119.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 81765.38
$ws.Range("I70").Value = 158223.08
$ws.Range("J70").Value = 5307.6924
$ws.Range("K70").Value = 158223.08
$ws.Range("L70").Value = 5307.6924
$ws.Range("M70").Value = -157953.08
$ws.Range("N70").Value = -5847.6924

$ws.Range("H73").Value = 81765.38
$ws.Range("I73").Value = 158223.08
$ws.Range("J73").Value = 5307.6924
$ws.Range("K73").Value = 158223.08
$ws.Range("L73").Value = 5307.6924
$ws.Range("M73").Value = -157287.08
$ws.Range("N73").Value = -7179.6924

$ws.Range("H110").Value = 39200
$ws.Range("J110").Value = 39200
$ws.Range("L110").Value = 39200
$ws.Range("N110").Value = -47380

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2823.5264
$ws.Range("J68").Value = 3319
$ws.Range("L68").Value = 3319
$ws.Range("N68").Value = -4817

$ws.Range("H69").Value = 34950
$ws.Range("J69").Value = 34950
$ws.Range("L69").Value = 34950
$ws.Range("N69").Value = -36572

$ws.Range("H71").Value = 2823.5264
$ws.Range("J71").Value = 3319
$ws.Range("L71").Value = 16595
$ws.Range("N71").Value = -24083

$ws.Range("H72").Value = 34950
$ws.Range("J72").Value = 34950
$ws.Range("L72").Value = 104850
$ws.Range("N72").Value = -112962

$ws.Range("H131").Value = 26285.715
$ws.Range("J131").Value = 26285.715
$ws.Range("L131").Value = 26285.715
$ws.Range("N131").Value = -36365.715

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H26").Value = 500000
$ws.Range("J26").Value = 500000
$ws.Range("L26").Value = 500000
$ws.Range("N26").Value = -500586

$ws.Range("H42").Value = 70049
$ws.Range("J42").Value = 70049
$ws.Range("L42").Value = 70049
$ws.Range("N42").Value = -70805

$ws.Range("H62").Value = 31260000
$ws.Range("I62").Value = 62500000
$ws.Range("J62").Value = 20000
$ws.Range("K62").Value = 62500000
$ws.Range("L62").Value = 20000
$ws.Range("M62").Value = -62499376
$ws.Range("N62").Value = -21248

$ws.Range("H65").Value = 31260000
$ws.Range("I65").Value = 62500000
$ws.Range("J65").Value = 20000
$ws.Range("K65").Value = 312500000
$ws.Range("L65").Value = 100000
$ws.Range("M65").Value = -312496880
$ws.Range("N65").Value = -106240

